# "Generate Report for Handback"
#
# The handback process for file d5f8d9c5-fddc-4138-95cc-ced84f90c5ab (row 7
# in the per-language sheets) failed its transform, because the file name
# coming back from the translator didn't match the file name that was
# handed off. Update the Status column for that row (on the Overview sheet
# and both language sheets) and record the error detail in column K of the
# two language sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"

# Overview sheet: both the zh-cn and de-de status columns for this file
# show the same aggregate status.
$overview.Range("B7").Value = $newStatus
$overview.Range("C7").Value = $newStatus

# zh-cn sheet: update status and add the error detail explaining why the
# handback transform failed.
$zhcn.Range("C7").Value = $newStatus
$zhcn.Range("K7").Value = "Handback file name: yfet0u0f.mjh is different with handoff file name: d5f8d9c5-fddc-4138-95cc-ced84f90c5ab.de8834329904febd8a1e80b8d9eb4113543aa339.zh-cn."

# de-de sheet: same update for the German locale row.
$dede.Range("C7").Value = $newStatus
$dede.Range("K7").Value = "Handback file name: yfet0u0f.mjh is different with handoff file name: d5f8d9c5-fddc-4138-95cc-ced84f90c5ab.de8834329904febd8a1e80b8d9eb4113543aa339.de-de."
